$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correction: 15/04/2020 case count updated from 832 to 882
$ws.Range("B32").Value = 882

# Insert the new 16/04/2020 data point as a new row right after 16/03/2020 (row 33),
# pushing 17/03/2020..31/03/2020 down by one row.
$ws.Range("A34").EntireRow.Insert()
$ws.Range("A34").Value = "16/04/2020"
$ws.Range("B34").Value = 190
